# edit.ps1 - Apply the commit's changes via PowerPoint COM-interop:
#   1. Change the table style (tableStyleId) on the table in slide 16
#      from {F75630C4-4F52-4FF3-97E4-E59ABED1A619} to
#      {55C4F410-3B66-45B7-98CE-4DEF579317FC}.
#   2. Swap the two themes used by the deck: the theme actually driving
#      the slide master/slides ("Integral") is replaced with the colors
#      of the "Office Theme" palette that the deck also carries.

$p = $ppt.ActivePresentation

# --- 1. Table style update (slide 16, 3rd shape = the graphicFrame/table) ---
$slide = $p.Slides.Item(16)
$shape = $slide.Shapes.Item(3)
if ($shape.HasTable) {
    $table = $shape.Table
    $table.ApplyStyle("{55C4F410-3B66-45B7-98CE-4DEF579317FC}")
}

# --- 2. Theme color swap (Integral -> Office) on the slide master's theme ---
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

function HexToRGBVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return RGBVal $r $g $b
}

# Office Theme color scheme (currently sitting unused in theme1.xml);
# index order matches ThemeColorScheme.Item(1..12):
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToRGBVal $officeColors[$i - 1]
}
